$wb = $excel.ActiveWorkbook

# ---- ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 472.0435
$ws.Range("I33").Value = 407.5238
$ws.Range("J33").Value = 1149.5
$ws.Range("K33").Value = 407.5238
$ws.Range("L33").Value = 1149.5
$ws.Range("M33").Value = -178.5238
$ws.Range("N33").Value = -1607.5

# ---- ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 15251
$ws.Range("I2").Value = 20524.268
$ws.Range("K2").Value = 20524.268
$ws.Range("M2").Value = -20411.268
$ws.Range("H32").Value = 7745.8887
$ws.Range("I32").Value = 5671.7734
$ws.Range("K32").Value = 5671.7734
$ws.Range("M32").Value = -5384.7734
$ws.Range("H110").Value = 1215.3914
$ws.Range("I110").Value = 1262.2727
$ws.Range("K110").Value = 1262.2727
$ws.Range("M110").Value = 782.7273
$ws.Range("H116").Value = 15251
$ws.Range("I116").Value = 20524.268
$ws.Range("K116").Value = 20524.268
$ws.Range("M116").Value = -18230.268

# ---- BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 15251
$ws.Range("I3").Value = 20524.268
$ws.Range("K3").Value = 20524.268
$ws.Range("M3").Value = -20410.268
$ws.Range("H20").Value = 80615.30499999999
$ws.Range("I20").Value = 2379.625
$ws.Range("J20").Value = 205792.4
$ws.Range("K20").Value = 2379.625
$ws.Range("L20").Value = 205792.4
$ws.Range("M20").Value = -2132.625
$ws.Range("N20").Value = -206286.4
$ws.Range("H86").Value = 7807.9
$ws.Range("I86").Value = 2821.6667
$ws.Range("J86").Value = 9944.857
$ws.Range("K86").Value = 2821.6667
$ws.Range("L86").Value = 9944.857
$ws.Range("M86").Value = -1698.6667
$ws.Range("N86").Value = -12190.857
$ws.Range("H89").Value = 7807.9
$ws.Range("I89").Value = 2821.6667
$ws.Range("J89").Value = 9944.857
$ws.Range("K89").Value = 14108.3335
$ws.Range("L89").Value = 49724.285
$ws.Range("M89").Value = -8492.333500000001
$ws.Range("N89").Value = -60956.285
$ws.Range("H94").Value = 3134.52
$ws.Range("I94").Value = 2934.4
$ws.Range("J94").Value = 3935
$ws.Range("K94").Value = 2934.4
$ws.Range("L94").Value = 3935
$ws.Range("M94").Value = -2483.4
$ws.Range("N94").Value = -4837
$ws.Range("H106").Value = 30223.666
$ws.Range("J106").Value = 30223.666
$ws.Range("L106").Value = 30223.666
$ws.Range("N106").Value = -32747.666
$ws.Range("H107").Value = 732
$ws.Range("I107").Value = 732
$ws.Range("K107").Value = 732
$ws.Range("M107").Value = 1188
$ws.Range("H113").Value = 0
$ws.Range("I113").Value = 0
$ws.Range("K113").Value = 0
$ws.Range("M113").ClearContents()
$ws.Range("H134").Value = 4674.6
$ws.Range("I134").Value = 3624.4736
$ws.Range("K134").Value = 10873.4208
$ws.Range("M134").Value = -8338.4208

# ---- CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 3095.7368
$ws.Range("I132").Value = 3078.889
$ws.Range("K132").Value = 9236.667000000001
$ws.Range("M132").Value = -6706.667000000001
$ws.Range("H141").Value = 89448.27
$ws.Range("J141").Value = 89448.27
$ws.Range("L141").Value = 89448.27
$ws.Range("N141").Value = -99808.27

# ---- CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H17").Value = 9999
$ws.Range("I17").Value = 9999
$ws.Range("J17").Value = 0
$ws.Range("K17").Value = 29997
$ws.Range("L17").Value = 0
$ws.Range("M17").Value = -29828
$ws.Range("N17").ClearContents()
$ws.Range("H34").Value = 2230.8333
$ws.Range("I34").Value = 346.25
$ws.Range("J34").Value = 6000
$ws.Range("K34").Value = 1038.75
$ws.Range("L34").Value = 18000
$ws.Range("M34").Value = -954.75
$ws.Range("N34").Value = -18168
$ws.Range("H39").Value = 4165.4287
$ws.Range("J39").Value = 4165.4287
$ws.Range("L39").Value = 12496.2861
$ws.Range("N39").Value = -13084.2861
$ws.Range("H55").Value = 1003551.3
$ws.Range("J55").Value = 1253938.6
$ws.Range("L55").Value = 3761815.8
$ws.Range("N55").Value = -3762169.8
$ws.Range("H57").Value = 4450
$ws.Range("J57").Value = 4450
$ws.Range("L57").Value = 13350
$ws.Range("N57").Value = -14468
$ws.Range("H60").Value = 326.4074
$ws.Range("I60").Value = 350.33334
$ws.Range("K60").Value = 1051.00002
$ws.Range("M60").Value = -800.0000199999999
$ws.Range("H132").Value = 1119.3334
$ws.Range("I132").Value = 1053.8572
$ws.Range("J132").Value = 1348.5
$ws.Range("K132").Value = 9484.7148
$ws.Range("L132").Value = 12136.5
$ws.Range("M132").Value = -6954.7148
$ws.Range("N132").Value = -17196.5
$ws.Range("H139").Value = 2644.2
$ws.Range("I139").Value = 2604.6667
$ws.Range("J139").Value = 3000
$ws.Range("K139").Value = 7814.000100000001
$ws.Range("L139").Value = 9000
$ws.Range("M139").Value = -2674.000100000001
$ws.Range("N139").Value = -19280

# ---- GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 89862.30499999999
$ws.Range("I70").Value = 162730.86
$ws.Range("K70").Value = 162730.86
$ws.Range("M70").Value = -162460.86
$ws.Range("H73").Value = 89862.30499999999
$ws.Range("I73").Value = 162730.86
$ws.Range("K73").Value = 162730.86
$ws.Range("M73").Value = -161794.86
$ws.Range("H126").Value = 2546.818
$ws.Range("I126").Value = 2091.5557
$ws.Range("J126").Value = 4595.5
$ws.Range("K126").Value = 6274.6671
$ws.Range("L126").Value = 13786.5
$ws.Range("M126").Value = -3804.6671
$ws.Range("N126").Value = -18726.5

# ---- LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 5928.037
$ws.Range("I7").Value = 6595.5713
$ws.Range("K7").Value = 6595.5713
$ws.Range("M7").Value = -6483.5713
$ws.Range("H36").Value = 72539.664
$ws.Range("J36").Value = 72539.664
$ws.Range("L36").Value = 72539.664
$ws.Range("N36").Value = -73663.664
$ws.Range("H40").Value = 11331.737
$ws.Range("I40").Value = 15519.125
$ws.Range("J40").Value = 8286.362999999999
$ws.Range("K40").Value = 15519.125
$ws.Range("L40").Value = 8286.362999999999
$ws.Range("M40").Value = -15383.125
$ws.Range("N40").Value = -8558.362999999999
$ws.Range("H68").Value = 6360.2856
$ws.Range("J68").Value = 6239.4
$ws.Range("L68").Value = 6239.4
$ws.Range("N68").Value = -7737.4
$ws.Range("H71").Value = 6360.2856
$ws.Range("J71").Value = 6239.4
$ws.Range("L71").Value = 31197
$ws.Range("N71").Value = -38685
$ws.Range("H124").Value = 265000
$ws.Range("J124").Value = 265000
$ws.Range("L124").Value = 265000
$ws.Range("N124").Value = -274820
$ws.Range("H125").Value = 69999.8
$ws.Range("J125").Value = 69999.8
$ws.Range("L125").Value = 69999.8
$ws.Range("N125").Value = -79839.8
$ws.Range("H126").Value = 5928.037
$ws.Range("I126").Value = 6595.5713
$ws.Range("K126").Value = 19786.7139
$ws.Range("M126").Value = -17316.7139
$ws.Range("H127").Value = 0
$ws.Range("J127").Value = 0
$ws.Range("L127").Value = 0
$ws.Range("N127").ClearContents()

# ---- WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 1524.7142
$ws.Range("J62").Value = 1339.6
$ws.Range("L62").Value = 1339.6
$ws.Range("N62").Value = -2587.6
$ws.Range("H65").Value = 1524.7142
$ws.Range("J65").Value = 1339.6
$ws.Range("L65").Value = 6698
$ws.Range("N65").Value = -12938
$ws.Range("H107").Value = 484.5
$ws.Range("I107").Value = 522
$ws.Range("K107").Value = 1566
$ws.Range("M107").Value = 354
$ws.Range("H126").Value = 1828.1818
$ws.Range("I126").Value = 1702.4736
$ws.Range("K126").Value = 5107.4208
$ws.Range("M126").Value = -2637.4208
